# The commit message says "removed spaces and quotation marks from files" -
# the three header rows of the VO2 export (rows 1-3) had labels/units padded
# with literal spaces and wrapped in escaped double-quotes (e.g. ' "VO2   "',
# 'TIME ', 'min  '). This script rewrites those header cells to their clean,
# trimmed text, clears the cells that were only quote/space filler (so they
# disappear from the sheet, matching the diff), and also fixes the two
# previously-identical "VE/" labels (M1/N1) to the distinct "VE/VO2" /
# "VE/VCO2" they should have been. Row 4's "----------" separator is
# untouched. Finally the active selection is moved to N2, matching the
# saved view state in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column headers -------------------------------------------------
$ws.Range("A1").Value = "TIME"
$ws.Range("B1").Value = "VO2"
$ws.Range("C1").Value = "VO2/kg"
$ws.Range("D1").Value = "METS"
$ws.Range("E1").Value = "VCO2"
$ws.Range("F1").Value = "VE"
$ws.Range("G1").Value = "RER"
$ws.Range("H1").Value = "RR"
$ws.Range("I1").Value = "Vt"
$ws.Range("J1").Value = "FEO2"
$ws.Range("K1").Value = "FECO2"
$ws.Range("L1").Value = "HR"
$ws.Range("M1").Value = "VE/VO2"
$ws.Range("N1").Value = "VE/VCO2"
$ws.Range("O1").Value = "PetCO2"
$ws.Range("P1").Value = "PO"

# --- Row 2: secondary unit/condition row ------------------------------------
# Cells that held nothing but spaces/quotes are cleared outright.
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "STPD"
$ws.Range("C2").Value = "STPD"
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "STPD"
$ws.Range("F2").Value = "BTPS"
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = "BTPS"
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()

# --- Row 3: units row -------------------------------------------------------
$ws.Range("A3").Value = "min"
$ws.Range("B3").Value = "L/min"
$ws.Range("C3").Value = "ml/kg/m"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "L/min"
$ws.Range("F3").Value = "L/min"
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = "BPM"
$ws.Range("I3").Value = "L"
$ws.Range("J3").Value = "%"
$ws.Range("K3").Value = "%"
$ws.Range("L3").Value = "bpm"
$ws.Range("M3").Value = "BT/ST"
$ws.Range("N3").Value = "BT/ST"
$ws.Range("O3").Value = "mmHg"
$ws.Range("P3").Value = "W"
$ws.Range("Q3").Value = "Lactate"

# --- Row 4: separator (unchanged, left as-is) -------------------------------
$ws.Range("A4").Value = "----------"

# --- Final selection, matching the saved view state -------------------------
$ws.Range("N2").Select()
